$d = $word.ActiveDocument

# The last paragraph in the body is "Add MoreWillieHand and make it
# partial class WillieHand." -- append a brand-new list paragraph right
# after it containing "Exercise 2 - Part B". Splitting the paragraph mark
# via InsertParagraphAfter naturally carries over the same ListParagraph
# style / numbered-list (ilvl 0, numId 1) formatting, exactly like typing
# Enter at the end of that bullet in Word would.

$lastPara = $d.Paragraphs.Last
$insertRange = $lastPara.Range
$insertRange.Collapse(0)          # wdCollapseEnd

$insertRange.InsertParagraphAfter()
$insertRange.Collapse(0)          # wdCollapseEnd, move past the new paragraph mark

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Exercise 2 - Part B"
